$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all touched D/E cells to Text format first so values like "1.010" or
# "26.978.24" are stored as literal strings instead of being reinterpreted as
# numbers/dates by Excel, then clear the temporary format so the cell style
# index matches the original (unstyled) cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.978.24'
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").Value = '1.820.58'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").Value = '309.92'
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '0.4636'
$ws.Range("E7").Value = '  -2.70%  '
$ws.Range("D8").Value = '0.3640'
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").Value = '0.07290'
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").Value = '0.8665'
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("D11").Value = '19.82'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("D12").Value = '1.883.43'
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").Value = '0.07619'
$ws.Range("E13").Value = '  +3.39%  '
$ws.Range("D14").Value = '93.17'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '5.331'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").Value = '6.479'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = '0.000008629'
$ws.Range("E18").Value = '  -2.55%  '
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").Value = '27.364.38'
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = '5.161'
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").Value = '2.110.16'
$ws.Range("E24").Value = '  +1.69%  '
$ws.Range("D25").Value = '151.87'
$ws.Range("D26").Value = '1.860'
$ws.Range("E26").Value = '  -2.10%  '
$ws.Range("E27").Value = '  -2.28%  '
$ws.Range("D28").Value = '2.098'
$ws.Range("E28").Value = '  -3.40%  '
$ws.Range("D29").Value = '5.090'
$ws.Range("E29").Value = '  -3.57%  '
$ws.Range("D30").Value = '115.91'
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("D31").Value = '0.08905'
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = '2.952'
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").Value = '0.7305'
$ws.Range("E33").Value = '  -3.96%  '
$ws.Range("D34").Value = '1.142'
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("D35").Value = '4.433'
$ws.Range("E35").Value = '  -2.87%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  +5.17%  '
$ws.Range("D38").Value = '0.05278'
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").Value = '0.01917'
$ws.Range("E40").Value = '  -2.51%  '
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("D42").Value = '7.161'
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").Value = '0.5222'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = '0.1634'
$ws.Range("E44").Value = '  -2.10%  '
$ws.Range("D45").Value = '8.257'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("D46").Value = '0.4862'
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("D47").Value = '1.009'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").Value = '10.16'
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").Value = '103.29'
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("D50").Value = '1.636'
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("D51").Value = '0.06222'
$ws.Range("E51").Value = '  -1.64%  '

# Remove the temporary text-number-format so the cells end up with the same
# (absence of) style index as before the edit.
$ws.Range("D2:E51").ClearFormats()
